# Scheduled market-data refresh for Kraken_Profits sheets.
# Updates computed price/profit columns (H:N) per leve row across all
# craft-class worksheets. Values come from the latest market snapshot;
# when a column has no meaningful value for a row it is cleared rather
# than left at a stale figure.

$wb = $excel.ActiveWorkbook

# ==== Sheet: ALC ====
$ws = $wb.Worksheets.Item("ALC")

# Row 5
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = ""
$ws.Range("N5").Value = ""

# Row 6
$ws.Range("H6").Value = 15.5
$ws.Range("I6").Value = 15.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 46.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 65.5

# Row 17
$ws.Range("H17").Value = 1619.7273
$ws.Range("I17").Value = 1549
$ws.Range("J17").Value = 1626.8
$ws.Range("K17").Value = 4647
$ws.Range("L17").Value = 4880.4
$ws.Range("M17").Value = -4479
$ws.Range("N17").Value = -5216.4

# Row 86
$ws.Range("H86").Value = 2974.25
$ws.Range("I86").Value = 2999
$ws.Range("J86").Value = 2900
$ws.Range("K86").Value = 2999
$ws.Range("L86").Value = 2900
$ws.Range("M86").Value = -1876
$ws.Range("N86").Value = -5146

# Row 88
$ws.Range("H88").Value = 1586.8334
$ws.Range("I88").Value = 1913.75
$ws.Range("J88").Value = 933
$ws.Range("K88").Value = 1913.75
$ws.Range("L88").Value = 933
$ws.Range("M88").Value = -1507.75
$ws.Range("N88").Value = -1745

# Row 89
$ws.Range("H89").Value = 2974.25
$ws.Range("I89").Value = 2999
$ws.Range("J89").Value = 2900
$ws.Range("K89").Value = 14995
$ws.Range("L89").Value = 14500
$ws.Range("M89").Value = -9379
$ws.Range("N89").Value = -25732

# Row 91
$ws.Range("H91").Value = 1586.8334
$ws.Range("I91").Value = 1913.75
$ws.Range("J91").Value = 933
$ws.Range("K91").Value = 1913.75
$ws.Range("L91").Value = 933
$ws.Range("M91").Value = -509.75
$ws.Range("N91").Value = -3741

# Row 92
$ws.Range("H92").Value = 4123.5
$ws.Range("I92").Value = 1500
$ws.Range("J92").Value = 4998
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 4998
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -7494

# Row 115
$ws.Range("H115").Value = 485
$ws.Range("I115").Value = 485
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 1455
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = 112

# Row 116
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = ""

# Row 125
$ws.Range("H125").Value = 852.1111
$ws.Range("I125").Value = 727.6
$ws.Range("J125").Value = 1007.75
$ws.Range("K125").Value = 6548.400000000001
$ws.Range("L125").Value = 9069.75
$ws.Range("M125").Value = -4088.400000000001
$ws.Range("N125").Value = -13989.75

# Row 132
$ws.Range("H132").Value = 5380.409
$ws.Range("I132").Value = 4658.8667
$ws.Range("J132").Value = 6926.5713
$ws.Range("K132").Value = 13976.6001
$ws.Range("L132").Value = 20779.7139
$ws.Range("M132").Value = -11446.6001
$ws.Range("N132").Value = -25839.7139

# Row 137
$ws.Range("H137").Value = 6499
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 6499
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 19497
$ws.Range("M137").Value = ""

# ==== Sheet: ARM ====
$ws = $wb.Worksheets.Item("ARM")

# Row 11
$ws.Range("H11").Value = 6004599.5
$ws.Range("I11").Value = 30000000
$ws.Range("J11").Value = 5749.5
$ws.Range("K11").Value = 30000000
$ws.Range("L11").Value = 5749.5
$ws.Range("M11").Value = -29999856
$ws.Range("N11").Value = -6037.5

# Row 61
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = ""
$ws.Range("N61").Value = ""

# Row 92
$ws.Range("H92").Value = 35000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 35000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 35000
$ws.Range("N92").Value = -39992

# Row 132
$ws.Range("H132").Value = 4273.125
$ws.Range("I132").Value = 2683.75
$ws.Range("J132").Value = 5862.5
$ws.Range("K132").Value = 8051.25
$ws.Range("L132").Value = 17587.5
$ws.Range("M132").Value = -5521.25
$ws.Range("N132").Value = -22647.5

# Row 136
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = ""
$ws.Range("N136").Value = ""

# ==== Sheet: BSM ====
$ws = $wb.Worksheets.Item("BSM")

# Row 20
$ws.Range("H20").Value = 1200
$ws.Range("I20").Value = 500
$ws.Range("J20").Value = 1433.3334
$ws.Range("K20").Value = 500
$ws.Range("L20").Value = 1433.3334
$ws.Range("M20").Value = -253

# Row 108
$ws.Range("H108").Value = 40000
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 40000
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 40000
$ws.Range("N108").Value = -47680

# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = ""

# Row 134
$ws.Range("H134").Value = 9666
$ws.Range("I134").Value = 8999.5
$ws.Range("J134").Value = 10999
$ws.Range("K134").Value = 26998.5
$ws.Range("L134").Value = 32997
$ws.Range("M134").Value = -24463.5
$ws.Range("N134").Value = -38067

# ==== Sheet: CRP ====
$ws = $wb.Worksheets.Item("CRP")

# Row 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").Value = ""

# Row 41
$ws.Range("H41").Value = 13766.333
$ws.Range("I41").Value = 1299
$ws.Range("J41").Value = 20000
$ws.Range("K41").Value = 1299
$ws.Range("L41").Value = 20000
$ws.Range("M41").Value = -871
$ws.Range("N41").Value = -20856

# Row 134
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = ""
$ws.Range("N134").Value = ""

# ==== Sheet: CUL ====
$ws = $wb.Worksheets.Item("CUL")

# Row 11
$ws.Range("H11").Value = 7143260.5
$ws.Range("I11").Value = 10000380
$ws.Range("J11").Value = 462.5
$ws.Range("K11").Value = 30001140
$ws.Range("L11").Value = 1387.5
$ws.Range("M11").Value = -30001000

# Row 26
$ws.Range("H26").Value = 540.2
$ws.Range("I26").Value = 416.33334
$ws.Range("J26").Value = 726
$ws.Range("K26").Value = 1249.00002
$ws.Range("L26").Value = 2178
$ws.Range("M26").Value = -961.0000199999999
$ws.Range("N26").Value = -2754

# Row 81
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = ""

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = ""

# Row 109
$ws.Range("H109").Value = 227
$ws.Range("I109").Value = 227
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 681
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = 359

# Row 115
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = ""

# Row 131
$ws.Range("H131").Value = 1674.7142
$ws.Range("I131").Value = 1159.2858
$ws.Range("J131").Value = 2190.1428
$ws.Range("K131").Value = 3477.8574
$ws.Range("L131").Value = 6570.428400000001
$ws.Range("M131").Value = 1562.1426
$ws.Range("N131").Value = -16650.4284

# Row 132
$ws.Range("H132").Value = 648.3333
$ws.Range("I132").Value = 648.3333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5834.9997
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3304.9997
$ws.Range("N132").Value = ""

# Row 134
$ws.Range("H134").Value = 1600
$ws.Range("I134").Value = 1600
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4800
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 270

# Row 139
$ws.Range("H139").Value = 5116.25
$ws.Range("I139").Value = 930
$ws.Range("J139").Value = 5714.2856
$ws.Range("K139").Value = 2790
$ws.Range("L139").Value = 17142.8568
$ws.Range("M139").Value = 2350
$ws.Range("N139").Value = -27422.8568

# Row 140
$ws.Range("H140").Value = 811
$ws.Range("I140").Value = 811
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 2433
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 2747

# ==== Sheet: GSM ====
$ws = $wb.Worksheets.Item("GSM")

# Row 70
$ws.Range("H70").Value = 4749.5
$ws.Range("I70").Value = 4499.5
$ws.Range("J70").Value = 4999.5
$ws.Range("K70").Value = 4499.5
$ws.Range("L70").Value = 4999.5
$ws.Range("M70").Value = -4229.5

# Row 73
$ws.Range("H73").Value = 4749.5
$ws.Range("I73").Value = 4499.5
$ws.Range("J73").Value = 4999.5
$ws.Range("K73").Value = 4499.5
$ws.Range("L73").Value = 4999.5
$ws.Range("M73").Value = -3563.5

# Row 80
$ws.Range("H80").Value = 3255.5
$ws.Range("I80").Value = 3255.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3255.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2257.5

# Row 83
$ws.Range("H83").Value = 3255.5
$ws.Range("I83").Value = 3255.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 16277.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -11285.5

# Row 126
$ws.Range("H126").Value = 3262.3333
$ws.Range("I126").Value = 3899.75
$ws.Range("J126").Value = 1987.5
$ws.Range("K126").Value = 11699.25
$ws.Range("L126").Value = 5962.5
$ws.Range("M126").Value = -9229.25
$ws.Range("N126").Value = -10902.5

# Row 132
$ws.Range("H132").Value = 4436.091
$ws.Range("I132").Value = 4144.222
$ws.Range("J132").Value = 5749.5
$ws.Range("K132").Value = 12432.666
$ws.Range("L132").Value = 17248.5
$ws.Range("M132").Value = -9902.665999999999
$ws.Range("N132").Value = -22308.5

# ==== Sheet: LTW ====
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 4666.6665
$ws.Range("I22").Value = 2000
$ws.Range("J22").Value = 6000
$ws.Range("K22").Value = 2000
$ws.Range("L22").Value = 6000
$ws.Range("M22").Value = -1705
$ws.Range("N22").Value = -6590

# Row 27
$ws.Range("H27").Value = 4666.6665
$ws.Range("I27").Value = 2000
$ws.Range("J27").Value = 6000
$ws.Range("K27").Value = 2000
$ws.Range("L27").Value = 6000
$ws.Range("M27").Value = -1893
$ws.Range("N27").Value = -6214

# ==== Sheet: WVR ====
$ws = $wb.Worksheets.Item("WVR")

# Row 10
$ws.Range("H10").Value = 5001002.5
$ws.Range("I10").Value = 5001002.5
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 5001002.5
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -5000833.5

# Row 45
$ws.Range("H45").Value = 42745.668
$ws.Range("I45").Value = 33250
$ws.Range("J45").Value = 47493.5
$ws.Range("K45").Value = 33250
$ws.Range("L45").Value = 47493.5
$ws.Range("M45").Value = -32759
$ws.Range("N45").Value = -48475.5

# Row 104
$ws.Range("H104").Value = 4990
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 4990
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 4990
$ws.Range("N104").Value = -11978
